$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.432.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "195.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.630"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.77"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000286"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.22"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.093.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "665.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +11.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.432.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.93%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.528.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.964"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.32"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.37"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.09"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.75"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.789.98"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.13%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.23%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.71"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "499.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.370"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.07%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.81%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.46%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.64%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.77"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +20.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +64.80%  "
